# Applies the "Windows Running.docx" edit:
#  1. Paragraph 1 ("Command line for running on windows") gets a
#     first-line indent and the _GoBack bookmark is relocated here
#     (it previously sat at the end of the "#include <winsock2.h>"
#     paragraph).
#  2. The _GoBack bookmark is removed from the "#include <winsock2.h>"
#     paragraph.
#  3. New "Combined:" section with two copies of the combine.c gcc
#     command line is appended at the end of the document.

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$d = $word.ActiveDocument

# --- 1. First paragraph: add first-line indent + move the bookmark in ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$xml1 = @"
<w:p xmlns:w="$wNs"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Command line for running on windows</w:t></w:r></w:p>
"@
$r1.InsertXML($xml1)

# --- 2. "#include <winsock2.h>" paragraph: drop the bookmark ---
$p8 = $d.Paragraphs.Item(8)
$r8 = $p8.Range
$xml8 = @"
<w:p xmlns:w="$wNs"><w:r><w:t>#include &lt;winsock2.h&gt;</w:t></w:r></w:p>
"@
$r8.InsertXML($xml8)

# --- 3. Append the new "Combined:" block after the trailing blank paragraph ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $pLast.Range
$gccCombined = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>gcc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>std</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">=c99 -o comb </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>combine.c</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -lws2_32 -lgdi32 -lcomdlg32</w:t></w:r>'
$xmlLast = @"
<w:p xmlns:w="$wNs"/><w:p xmlns:w="$wNs"><w:r><w:t>Combined:</w:t></w:r></w:p><w:p xmlns:w="$wNs">$gccCombined</w:p><w:p xmlns:w="$wNs"/><w:p xmlns:w="$wNs">$gccCombined</w:p>
"@
$rLast.InsertXML($xmlLast)
